# end of Rev 1 revision -- sec 5 plus nitpicking elsewhere
#
# Three small text nitpicks on slide 1 (everything else in the original
# commit -- the "2/24/20" -> "2/27/20" datetimeFigureOut field caches --
# lives in slide footers that this deck doesn't have, so there is nothing
# to touch for that part here):
#
#   - "TextBox 40"  (shape reading "g3y")  -> "3y" becomes "3z"
#   - "TextBox 53"  (shape reading "g2N")  -> the "2" is retyped as "3",
#                    so the shape now reads "g3N"
#   - "TextBox 70"  (shape reading "g1z")  -> "1z" becomes "1y"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "3y" -> "3z" --------------------------------------------------------
$shape1 = $s.Shapes.Item("TextBox 40")

# These little labels auto-fit to their text (spAutoFit); re-measuring the
# text after an edit can nudge the cached shape size by a hair, so stash
# the current height/width and restore them once the text is updated --
# the source diff only touches the run text, not the shape geometry.
$h1 = $shape1.Height
$w1 = $shape1.Width
$shape1.TextFrame.TextRange.Runs(2).Text = "3z"
$shape1.Height = $h1
$shape1.Width = $w1

# --- "2N" -> "3" + "N" (the leading "2" is retyped as "3") --------------
$shape2 = $s.Shapes.Item("TextBox 53")
$shape2.TextFrame.TextRange.Characters(2, 1).Text = "3"

# --- "1z" -> "1y" --------------------------------------------------------
$shape3 = $s.Shapes.Item("TextBox 70")
$shape3.TextFrame.TextRange.Runs(2).Text = "1y"
